{"js": "// Fix some running issues: stamp the (empty) first paragraph's run\n// properties with an explicit language of English (United States),\n// i.e. produce <w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>\n// on the lone paragraph in the document body.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The document body contains a single (empty) paragraph; set the\n// language on its range so Word records it on the paragraph mark's\n// run properties (w:rPr/w:lang).\nconst firstParagraph = paragraphs.items[0];\nconst range = firstParagraph.getRange();\nrange.languageId = \"en-US\";\n\nawait context.sync();\n", "ps1": "# Fix some running issues: stamp the (empty) first paragraph's run\n# properties with an explicit language of English (United States),\n# i.e. produce <w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>\n# on the lone paragraph in the document body.\n\n$d = $word.ActiveDocument\n\n# The document body contains a single (empty) paragraph; set the\n# language on its range so Word records it on the paragraph mark's\n# run properties (w:rPr/w:lang).\n$p = $d.Paragraphs(1)\n$r = $p.Range\n$r.LanguageID = \"en-US\"\n"}
